# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $text) {
    $c = $ws.Range($ref)
    # Leading apostrophe forces Excel to store the value as literal text
    # (prevents '437.80', '1.00', '0.0000142', etc. from being normalised as numbers),
    # then resetting the style back to Normal drops the transient quote-prefix flag
    # so the cell keeps its original (unstyled) appearance.
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '68.936.91'
Set-TextCell $ws 'E2' '  +2.22%  '
Set-TextCell $ws 'D3' '3.745.34'
Set-TextCell $ws 'E3' '  -1.10%  '
Set-TextCell $ws 'E4' '  +0.15%  '
Set-TextCell $ws 'D5' '601.27'
Set-TextCell $ws 'E5' '  +1.62%  '
Set-TextCell $ws 'D6' '168.42'
Set-TextCell $ws 'E6' '  -1.85%  '
Set-TextCell $ws 'D7' '3.745.66'
Set-TextCell $ws 'E7' '  -0.99%  '
Set-TextCell $ws 'E8' '  -0.01%  '
Set-TextCell $ws 'D9' '0.533'
Set-TextCell $ws 'E9' '  +2.22%  '
Set-TextCell $ws 'E10' '  +4.31%  '
Set-TextCell $ws 'E11' '  +1.28%  '
Set-TextCell $ws 'D12' '0.461'
Set-TextCell $ws 'E12' '  +0.12%  '
Set-TextCell $ws 'D13' '38.33'
Set-TextCell $ws 'E13' '  +1.77%  '
Set-TextCell $ws 'D14' '0.0000245'
Set-TextCell $ws 'E14' '  +0.49%  '
Set-TextCell $ws 'D15' '4.370.88'
Set-TextCell $ws 'E15' '  -1.02%  '
Set-TextCell $ws 'D16' '3.741.43'
Set-TextCell $ws 'E16' '  -0.86%  '
Set-TextCell $ws 'D17' '68.960.57'
Set-TextCell $ws 'E17' '  +2.15%  '
Set-TextCell $ws 'D18' '7.28'
Set-TextCell $ws 'E18' '  +2.18%  '
Set-TextCell $ws 'E19' '  +0.39%  '
Set-TextCell $ws 'D20' '17.22'
Set-TextCell $ws 'E20' '  +7.61%  '
Set-TextCell $ws 'D21' '498.50'
Set-TextCell $ws 'E21' '  +2.48%  '
Set-TextCell $ws 'D22' '9.73'
Set-TextCell $ws 'E22' '  +6.17%  '
Set-TextCell $ws 'D23' '0.725'
Set-TextCell $ws 'E23' '  +0.49%  '
Set-TextCell $ws 'D24' '84.85'
Set-TextCell $ws 'E24' '  +1.10%  '
Set-TextCell $ws 'B25' 'PEPE'
Set-TextCell $ws 'C25' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell $ws 'D25' '0.0000142'
Set-TextCell $ws 'E25' '  +1.62%  '
Set-TextCell $ws 'B26' 'Fetch.AI'
Set-TextCell $ws 'C26' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D26' '2.31'
Set-TextCell $ws 'E26' '  -2.19%  '
Set-TextCell $ws 'D27' '12.29'
Set-TextCell $ws 'E27' '  +1.13%  '
Set-TextCell $ws 'D28' '10.12'
Set-TextCell $ws 'E28' '  -0.41%  '
Set-TextCell $ws 'E29' '  +0.06%  '
Set-TextCell $ws 'E30' '  +1.41%  '
Set-TextCell $ws 'E31' '  +1.46%  '
Set-TextCell $ws 'E32' '  +2.16%  '
Set-TextCell $ws 'D33' '31.78'
Set-TextCell $ws 'E33' '  -2.22%  '
Set-TextCell $ws 'D34' '3.885.48'
Set-TextCell $ws 'E34' '  -0.98%  '
Set-TextCell $ws 'D35' '0.109'
Set-TextCell $ws 'E35' '  +0.93%  '
Set-TextCell $ws 'D36' '3.673.73'
Set-TextCell $ws 'E36' '  -1.38%  '
Set-TextCell $ws 'D37' '1.00'
Set-TextCell $ws 'E37' '  +0.28%  '
Set-TextCell $ws 'E38' '  +0.44%  '
Set-TextCell $ws 'D39' '5.79'
Set-TextCell $ws 'E39' '  +1.22%  '
Set-TextCell $ws 'E40' '  -1.28%  '
Set-TextCell $ws 'E41' '  +0.70%  '
Set-TextCell $ws 'D42' '437.80'
Set-TextCell $ws 'E42' '  -2.92%  '
Set-TextCell $ws 'D43' '49.01'
Set-TextCell $ws 'E43' '  +0.38%  '
Set-TextCell $ws 'E44' '  -0.51%  '
Set-TextCell $ws 'D45' '2.89'
Set-TextCell $ws 'E45' '  +2.09%  '
Set-TextCell $ws 'D46' '8.40'
Set-TextCell $ws 'E46' '  +1.92%  '
Set-TextCell $ws 'E47' '  +0.02%  '
Set-TextCell $ws 'D48' '40.53'
Set-TextCell $ws 'E48' '  -2.29%  '
Set-TextCell $ws 'D49' '143.59'
Set-TextCell $ws 'E49' '  +1.91%  '
Set-TextCell $ws 'E50' '  +1.14%  '
Set-TextCell $ws 'D51' '2.750.39'
Set-TextCell $ws 'E51' '  -2.85%  '
